$wb = $excel.ActiveWorkbook

# ALC!row4
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1257.4166
$ws.Range("I4").Value = 337.8
$ws.Range("K4").Value = 337.8
$ws.Range("M4").Value = -223.8

# ALC!row38
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1946
$ws.Range("I38").Value = 335.4
$ws.Range("K38").Value = 1006.2
$ws.Range("M38").Value = -634.1999999999999

# ALC!row62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3197.32
$ws.Range("J62").Value = 7830
$ws.Range("L62").Value = 7830
$ws.Range("N62").Value = -9078

# ALC!row65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 3197.32
$ws.Range("J65").Value = 7830
$ws.Range("L65").Value = 39150
$ws.Range("N65").Value = -45390

# ALC!row86
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 7672.5454
$ws.Range("J86").Value = 13876
$ws.Range("L86").Value = 13876
$ws.Range("N86").Value = -16122

# ALC!row89
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 7672.5454
$ws.Range("J89").Value = 13876
$ws.Range("L89").Value = 69380
$ws.Range("N89").Value = -80612

# ALC!row112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1987.1034
$ws.Range("J112").Value = 2004.8148
$ws.Range("L112").Value = 6014.4444
$ws.Range("N112").Value = -8230.4444

# ALC!row129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 3183.1875
$ws.Range("I129").Value = 1031.2727
$ws.Range("J129").Value = 7917.4
$ws.Range("K129").Value = 3093.8181
$ws.Range("L129").Value = 23752.2
$ws.Range("M129").Value = 1906.1819
$ws.Range("N129").Value = -33752.2

# ALC!row135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1774.4722
$ws.Range("I135").Value = 590.03125
$ws.Range("K135").Value = 5310.28125
$ws.Range("M135").Value = -2775.28125

# ALC!row137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 23813096
$ws.Range("I137").Value = 38465650
$ws.Range("K137").Value = 115396950
$ws.Range("M137").Value = -115394400

# ALC!row138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 6266.3477
$ws.Range("I138").Value = 3790.6428
$ws.Range("J138").Value = 10117.444
$ws.Range("K138").Value = 11371.9284
$ws.Range("L138").Value = 30352.332
$ws.Range("M138").Value = -6231.928400000001
$ws.Range("N138").Value = -40632.33199999999

# ALC!row140
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H140").Value = 99222
$ws.Range("J140").Value = 99222
$ws.Range("L140").Value = 99222
$ws.Range("N140").Value = -109582

# ARM!row32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6938.268
$ws.Range("I32").Value = 7283.846
$ws.Range("J32").Value = 199.5
$ws.Range("K32").Value = 7283.846
$ws.Range("L32").Value = 199.5
$ws.Range("M32").Value = -6996.846
$ws.Range("N32").Value = -773.5

# ARM!row38
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H38").Value = 20000
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 20000
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 20000
$ws.Range("M38").ClearContents()
$ws.Range("N38").Value = -20934

# ARM!row45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4263.077
$ws.Range("I45").Value = 2281.5
$ws.Range("K45").Value = 2281.5
$ws.Range("M45").Value = -1904.5

# ARM!row61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 26926302
$ws.Range("I61").Value = 41179704
$ws.Range("K61").Value = 41179704
$ws.Range("M61").Value = -41179492

# ARM!row74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2543.75
$ws.Range("I74").Value = 2660.25
$ws.Range("J74").Value = 2116.5833
$ws.Range("K74").Value = 2660.25
$ws.Range("L74").Value = 2116.5833
$ws.Range("M74").Value = -1786.25
$ws.Range("N74").Value = -3864.5833

# ARM!row77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2543.75
$ws.Range("I77").Value = 2660.25
$ws.Range("J77").Value = 2116.5833
$ws.Range("K77").Value = 13301.25
$ws.Range("L77").Value = 10582.9165
$ws.Range("M77").Value = -8933.25
$ws.Range("N77").Value = -19318.9165

# ARM!row124
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H124").Value = 19000
$ws.Range("J124").Value = 19000
$ws.Range("L124").Value = 19000
$ws.Range("N124").Value = -28820

# ARM!row132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4466.17
$ws.Range("I132").Value = 4376.3477
$ws.Range("K132").Value = 13129.0431
$ws.Range("M132").Value = -10599.0431

# ARM!row133
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H133").Value = 83999.5
$ws.Range("J133").Value = 74998
$ws.Range("L133").Value = 74998
$ws.Range("N133").Value = -80058

# ARM!row136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 26926302
$ws.Range("I136").Value = 41179704
$ws.Range("K136").Value = 123539112
$ws.Range("M136").Value = -123536562

# BSM!row22
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1960
$ws.Range("I22").Value = 981.1111
$ws.Range("K22").Value = 981.1111
$ws.Range("M22").Value = -808.1111

# BSM!row134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2075.838
$ws.Range("I134").Value = 1865.5151
$ws.Range("J134").Value = 3811
$ws.Range("K134").Value = 5596.5453
$ws.Range("L134").Value = 11433
$ws.Range("M134").Value = -3061.5453
$ws.Range("N134").Value = -16503

# CRP!row22
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 443.8
$ws.Range("I22").Value = 367.375
$ws.Range("J22").Value = 749.5
$ws.Range("K22").Value = 367.375
$ws.Range("L22").Value = 749.5
$ws.Range("M22").Value = -17.375
$ws.Range("N22").Value = -1449.5

# CRP!row31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 27781638
$ws.Range("I31").Value = 45457864
$ws.Range("J31").Value = 4708.357
$ws.Range("K31").Value = 45457864
$ws.Range("L31").Value = 4708.357
$ws.Range("M31").Value = -45457569
$ws.Range("N31").Value = -5298.357

# CRP!row34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 27781638
$ws.Range("I34").Value = 45457864
$ws.Range("J34").Value = 4708.357
$ws.Range("K34").Value = 45457864
$ws.Range("L34").Value = 4708.357
$ws.Range("M34").Value = -45457662
$ws.Range("N34").Value = -5112.357

# CRP!row99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 14942.286
$ws.Range("I99").Value = 8720.691999999999
$ws.Range("J99").Value = 25052.375
$ws.Range("K99").Value = 8720.691999999999
$ws.Range("L99").Value = 25052.375
$ws.Range("M99").Value = -7222.691999999999
$ws.Range("N99").Value = -28048.375

# CRP!row126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 14942.286
$ws.Range("I126").Value = 8720.691999999999
$ws.Range("J126").Value = 25052.375
$ws.Range("K126").Value = 26162.076
$ws.Range("L126").Value = 75157.125
$ws.Range("M126").Value = -23692.076
$ws.Range("N126").Value = -80097.125

# CRP!row132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1707.7435
$ws.Range("I132").Value = 972.7222
$ws.Range("K132").Value = 2918.1666
$ws.Range("M132").Value = -388.1666

# CRP!row134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1260.4147
$ws.Range("I134").Value = 1228.6578
$ws.Range("J134").Value = 1662.6666
$ws.Range("K134").Value = 3685.9734
$ws.Range("L134").Value = 4987.9998
$ws.Range("M134").Value = -1150.9734
$ws.Range("N134").Value = -10057.9998

# CUL!row2
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 79.888885
$ws.Range("I2").Value = 82
$ws.Range("J2").Value = 63
$ws.Range("K2").Value = 492
$ws.Range("L2").Value = 378
$ws.Range("M2").Value = -379
$ws.Range("N2").Value = -604

# CUL!row12
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 1202.7368
$ws.Range("J12").Value = 1569.9231
$ws.Range("L12").Value = 4709.7693
$ws.Range("N12").Value = -5055.7693

# CUL!row21
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H21").Value = 6321.7144
$ws.Range("J21").Value = 21666
$ws.Range("L21").Value = 64998
$ws.Range("N21").Value = -65344

# CUL!row92
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 170
$ws.Range("J92").Value = 181.66667
$ws.Range("L92").Value = 545.00001
$ws.Range("N92").Value = -3041.00001

# CUL!row120
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 14571.6
$ws.Range("I120").Value = 9881.25
$ws.Range("K120").Value = 29643.75
$ws.Range("M120").Value = -24805.75

# CUL!row132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 3508.889
$ws.Range("J132").Value = 4946
$ws.Range("L132").Value = 44514
$ws.Range("N132").Value = -49574

# GSM!row102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2322.84
$ws.Range("I102").Value = 2028.1177
$ws.Range("K102").Value = 2028.1177
$ws.Range("M102").Value = -406.1177

# LTW!row22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 12409.8
$ws.Range("I22").Value = 29137.25
$ws.Range("J22").Value = 1258.1666
$ws.Range("K22").Value = 29137.25
$ws.Range("L22").Value = 1258.1666
$ws.Range("M22").Value = -28842.25
$ws.Range("N22").Value = -1848.1666

# LTW!row27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 12409.8
$ws.Range("I27").Value = 29137.25
$ws.Range("J27").Value = 1258.1666
$ws.Range("K27").Value = 29137.25
$ws.Range("L27").Value = 1258.1666
$ws.Range("M27").Value = -29030.25
$ws.Range("N27").Value = -1472.1666

# LTW!row48
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H48").Value = 66681
$ws.Range("I48").Value = 99998
$ws.Range("K48").Value = 99998
$ws.Range("M48").Value = -99337

# LTW!row98
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H98").Value = 99997
$ws.Range("J98").Value = 99997
$ws.Range("L98").Value = 99997
$ws.Range("N98").Value = -105987

# LTW!row132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3351.258
$ws.Range("I132").Value = 1751.85
$ws.Range("J132").Value = 6259.273
$ws.Range("K132").Value = 5255.549999999999
$ws.Range("L132").Value = 18777.819
$ws.Range("M132").Value = -2725.549999999999
$ws.Range("N132").Value = -23837.819

# LTW!row136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2396.0667
$ws.Range("I136").Value = 1547.1428
$ws.Range("J136").Value = 4376.8887
$ws.Range("K136").Value = 4641.428400000001
$ws.Range("L136").Value = 13130.6661
$ws.Range("M136").Value = -2091.428400000001
$ws.Range("N136").Value = -18230.6661

# WVR!row132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2598.3333
$ws.Range("J132").Value = 3366.6667
$ws.Range("L132").Value = 10100.0001
$ws.Range("N132").Value = -15160.0001
